$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Delete Unit"

# Update the header cell text
$ws.Range("A1").Value = "Item Unit Name"

# Remove the now-unused second column ("Item Group Name")
$ws.Columns("B:B").Delete()
